# Weekly refresh of the Brócoli / Lo Valledor subset: a new week's worth of
# records (date 44783) is inserted at the top of this product's block, and
# all the existing rows for that block shift down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 796; this pushes the former rows 796-863 down to
# 798-865, carrying their values/styles with them (Excel's native Insert
# semantics), so we only need to populate the two brand-new rows below.
$ws.Rows.Item(796).Resize(2).Insert()

# New row 796: Brócoli, Primera, Región Metropolitana, fecha 44783 (2022-08-10)
$ws.Cells.Item(796, 1).Value = 6
$ws.Cells.Item(796, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(796, 3).Value = "Metropolitana"
$ws.Cells.Item(796, 4).Value = 44783
$ws.Cells.Item(796, 5).Value = 13
$ws.Cells.Item(796, 6).Value = 100112023
$ws.Cells.Item(796, 7).Value = "Brócoli"
$ws.Cells.Item(796, 8).Value = "Sin especificar"
$ws.Cells.Item(796, 9).Value = "Primera"
$ws.Cells.Item(796, 10).Value = 10600
$ws.Cells.Item(796, 11).Value = 500
$ws.Cells.Item(796, 12).Value = 600
$ws.Cells.Item(796, 13).Value = 546
$ws.Cells.Item(796, 14).Value = "`$/unidad"
$ws.Cells.Item(796, 15).Value = "Región Metropolitana"
$ws.Cells.Item(796, 16).Value = 546
$ws.Cells.Item(796, 17).Value = 1
$ws.Cells.Item(796, 18).Value = "Hortaliza"

# New row 797: Brócoli, Segunda, Región Metropolitana, fecha 44783 (2022-08-10)
$ws.Cells.Item(797, 1).Value = 6
$ws.Cells.Item(797, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(797, 3).Value = "Metropolitana"
$ws.Cells.Item(797, 4).Value = 44783
$ws.Cells.Item(797, 5).Value = 13
$ws.Cells.Item(797, 6).Value = 100112023
$ws.Cells.Item(797, 7).Value = "Brócoli"
$ws.Cells.Item(797, 8).Value = "Sin especificar"
$ws.Cells.Item(797, 9).Value = "Segunda"
$ws.Cells.Item(797, 10).Value = 3800
$ws.Cells.Item(797, 11).Value = 400
$ws.Cells.Item(797, 12).Value = 400
$ws.Cells.Item(797, 13).Value = 400
$ws.Cells.Item(797, 14).Value = "`$/unidad"
$ws.Cells.Item(797, 15).Value = "Región Metropolitana"
$ws.Cells.Item(797, 16).Value = 400
$ws.Cells.Item(797, 17).Value = 1
$ws.Cells.Item(797, 18).Value = "Hortaliza"
